$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# 1. Edit the task text in D5 first (so the shared string is mutated in place
#    and rows copied from it below pick up the new text under the same index).
$ws.Range("D5").Value = "Sanity testing on B2C app, QMVAR site, GSS site and Hayaai site. Regression testing, Retesting on Hayaai site and Sonia application"

# 2. Fill rows 6 and 7 by copying row 5 (carries over styles/number formats).
$ws.Range("A5:G5").Copy($ws.Range("A6:G6"))
$ws.Range("A5:G5").Copy($ws.Range("A7:G7"))

# 3. Update the incrementing sequence number and date for each new row.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 44140
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 44141

# 4. Match row height used by the rest of the filled-in rows.
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 30

# 5. Leave the active selection on the last edited cell.
$ws.Range("D7").Select()
